$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo "meeing" -> "meeting" in the cell currently holding that text (B12)
$ws.Range("B12").Value = "You are running short on time and don't understand why it is necessary to discuss this in private later. You could explain how you feel or just tell him the meeting won't take long."

# Update the active selection to B12 (was B13)
$ws.Range("B12").Select()
